$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "Changed fan angle late night feb 18. Returned fan to previous angle around noon feb 19."
$ws.Range("A45").Value = "Disconnected XXXL capacitor Sun Feb 19 around 7:20"
$ws.Range("A46").Value = "Opened window Feb 22 about 11:35AM"

[void]$ws.Range("A47").Select()
